# Apply the recorded edit to the "revenue&profit rate" workbook.
#
# Summary of the change (from the OOXML diff):
#  - Column M (rows 4-8) had a cell inserted at M4 with "shift down", so:
#      M4 (1052) -> blank
#      M5 (27236) -> 1052
#      M6 (1994) -> 27236
#      M7 (2306) -> 1994
#      M8 (31536) -> 2306
#  - Row 17 keeps its "B17" label cell but now also carries the (empty,
#    but styled) L17/M17 cells that used to live on the now-deleted row 18.
#  - Rows 18 and 19 (the stray leftover rows) are removed entirely, so the
#    sheet's used range shrinks from A1:N19 down to A1:N17.
#  - The view was re-zoomed (80%) and the selection moved to row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the M4:M8 values down by one cell (insert+shift-down on M4) ---
$ws.Range("M8").Value = $ws.Range("M7").Value()
$ws.Range("M7").Value = $ws.Range("M6").Value()
$ws.Range("M6").Value = $ws.Range("M5").Value()
$ws.Range("M5").Value = $ws.Range("M4").Value()
$ws.Range("M4").ClearContents()

# --- Give the new L17/M17 cells the same formatting as L7/M7 (style ids
#     "39"/"40" in the target file) by copying formats across, then
#     clearing any value so they stay blank. ---
$ws.Range("L7").Copy() | Out-Null
$ws.Range("L17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("M7").Copy() | Out-Null
$ws.Range("M17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("L17:M17").ClearContents()

# --- Remove the now-empty trailing rows 18 and 19 ---
$ws.Rows("18:19").Delete()

# --- Update the view: zoom to 80%, select row 7 ---
$excel.ActiveWindow.Zoom = 80
$ws.Range("A7:XFD7").Select() | Out-Null
